$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.204.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.824.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6019'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07116'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2809'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.05'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07649'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.876.82'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.766'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6398'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009685'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '79.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.039.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.988'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.136.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '229.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.72'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.017'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9998'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.067'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1278'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.63'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06779'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.454'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.458'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.810'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.759'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.130'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.716'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6567'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.535'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.766'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.78%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.223.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01758'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.521'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9259'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.961.67'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000116'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.627'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.557'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.99%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05580'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.485'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.59%  '
